# Fixed typo in learning_python 24 presentation
#
# 1) Update the cached "datetimeFigureOut" field text (12/10/2018 -> 18/12/2018)
#    on every slide layout that carries a Date Placeholder.
# 2) Fix the "de-bugging" typo (-> "debugging") and give "booleans" its own
#    run on slide 2's content placeholder - mirroring the way PowerPoint
#    splits a run when only part of a paragraph's text is retyped.

$p = $ppt.ActivePresentation

# --- 1. Date placeholders on the slide layouts ------------------------------
$master = $p.SlideMaster
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shape = $layout.Shapes.Item($j)
        if ($shape.HasTextFrame -and $shape.Name -like "Date Placeholder*") {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq "12/10/2018") {
                $tr.Text = "18/12/2018"
            }
        }
    }
}

# --- 2. Slide 2 bullet list typo fixes --------------------------------------
$slide2 = $p.Slides.Item(2)
for ($k = 1; $k -le $slide2.Shapes.Count; $k++) {
    $shape = $slide2.Shapes.Item($k)
    if ($shape.Name -eq "Content Placeholder 2") {
        $content = $shape
    }
}

$tr = $content.TextFrame.TextRange
$paraCount = $tr.Paragraphs().Count

for ($n = 1; $n -le $paraCount; $n++) {
    $para = $tr.Paragraphs($n, 1)
    # Paragraphs(...).Text includes the trailing paragraph mark (CR) - strip
    # it before comparing against the plain bullet text.
    $text = $para.Text.TrimEnd("`r", "`n")

    if ($text -eq "Basics and control flow, booleans") {
        $start = $text.IndexOf("booleans") + 1
        $len = "booleans".Length
        $para.Characters($start, $len).Text = "booleans"
    }
    elseif ($text -eq "Errors and de-bugging") {
        $start = $text.IndexOf("de-bugging") + 1
        $len = "de-bugging".Length
        $para.Characters($start, $len).Text = "debugging"
    }
}
